$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.572.63"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = "'  +5.79%  "
$ws.Range('E2').ClearFormats()
$ws.Range('D3').Value = "'1.723.67"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = "'  +4.40%  "
$ws.Range('E3').ClearFormats()
$ws.Range('E4').Value = "'  +0.09%  "
$ws.Range('E4').ClearFormats()
$ws.Range('D5').Value = "'225.96"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = "'  +3.51%  "
$ws.Range('E5').ClearFormats()
$ws.Range('D6').Value = "'0.5371"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = "'  +3.00%  "
$ws.Range('E6').ClearFormats()
$ws.Range('E7').Value = "'  +0.05%  "
$ws.Range('E7').ClearFormats()
$ws.Range('D8').Value = "'0.2669"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = "'  +1.16%  "
$ws.Range('E8').ClearFormats()
$ws.Range('D9').Value = "'0.06605"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = "'  +4.42%  "
$ws.Range('E9').ClearFormats()
$ws.Range('E10').Value = "'  +6.74%  "
$ws.Range('E10').ClearFormats()
$ws.Range('D11').Value = "'0.07720"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = "'  +0.73%  "
$ws.Range('E11').ClearFormats()
$ws.Range('D12').Value = "'4.624"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = "'  +0.80%  "
$ws.Range('E12').ClearFormats()
$ws.Range('B13').Value = "'WrappedliquidstakedEther2.0"
$ws.Range('B13').ClearFormats()
$ws.Range('C13').Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range('C13').ClearFormats()
$ws.Range('D13').Value = "'1.962.24"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = "'  +4.51%  "
$ws.Range('E13').ClearFormats()
$ws.Range('B14').Value = "'WrappedEther"
$ws.Range('B14').ClearFormats()
$ws.Range('C14').Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('C14').ClearFormats()
$ws.Range('D14').Value = "'1.702.54"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = "'  +2.13%  "
$ws.Range('E14').ClearFormats()
$ws.Range('D15').Value = "'0.5848"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = "'  +4.48%  "
$ws.Range('E15').ClearFormats()
$ws.Range('E16').Value = "'  +2.17%  "
$ws.Range('E16').ClearFormats()
$ws.Range('D17').Value = "'67.98"
$ws.Range('D17').ClearFormats()
$ws.Range('D18').Value = "'27.578.53"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = "'  +5.87%  "
$ws.Range('E18').ClearFormats()
$ws.Range('D19').Value = "'221.64"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = "'  +15.74%  "
$ws.Range('E19').ClearFormats()
$ws.Range('E20').Value = "'  +0.04%  "
$ws.Range('E20').ClearFormats()
$ws.Range('D21').Value = "'4.728"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = "'  +2.33%  "
$ws.Range('E21').ClearFormats()
$ws.Range('D22').Value = "'10.65"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = "'  +1.45%  "
$ws.Range('E22').ClearFormats()
$ws.Range('D23').Value = "'6.088"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = "'  +3.12%  "
$ws.Range('E23').ClearFormats()
$ws.Range('D24').Value = "'1.004"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = "'  +0.09%  "
$ws.Range('E24').ClearFormats()
$ws.Range('D25').Value = "'148.17"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = "'  +2.91%  "
$ws.Range('E25').ClearFormats()
$ws.Range('D26').Value = "'1.740"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = "'  +14.74%  "
$ws.Range('E26').ClearFormats()
$ws.Range('D27').Value = "'0.1236"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = "'  +4.23%  "
$ws.Range('E27').ClearFormats()
$ws.Range('D28').Value = "'7.417"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = "'  +3.11%  "
$ws.Range('E28').ClearFormats()
$ws.Range('D29').Value = "'16.61"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = "'  +4.67%  "
$ws.Range('E29').ClearFormats()
$ws.Range('D30').Value = "'0.05585"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = "'  +2.86%  "
$ws.Range('E30').ClearFormats()
$ws.Range('D31').Value = "'1.302"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = "'  +2.60%  "
$ws.Range('E31').ClearFormats()
$ws.Range('D32').Value = "'3.570"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = "'  +3.56%  "
$ws.Range('E32').ClearFormats()
$ws.Range('D33').Value = "'3.449"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = "'  +3.10%  "
$ws.Range('E33').ClearFormats()
$ws.Range('D34').Value = "'1.662"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = "'  +6.91%  "
$ws.Range('E34').ClearFormats()
$ws.Range('D35').Value = "'0.9692"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = "'  +2.58%  "
$ws.Range('E35').ClearFormats()
$ws.Range('D36').Value = "'2.840"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = "'  +2.12%  "
$ws.Range('E36').ClearFormats()
$ws.Range('D37').Value = "'2.426"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = "'  +0.15%  "
$ws.Range('E37').ClearFormats()
$ws.Range('D38').Value = "'0.5970"
$ws.Range('D38').ClearFormats()
$ws.Range('D39').Value = "'0.01654"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = "'  +4.63%  "
$ws.Range('E39').ClearFormats()
$ws.Range('D40').Value = "'5.909"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = "'  +0.87%  "
$ws.Range('E40').ClearFormats()
$ws.Range('D41').Value = "'0.8539"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = "'  +3.26%  "
$ws.Range('E41').ClearFormats()
$ws.Range('D42').Value = "'1.055.23"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = "'  +2.89%  "
$ws.Range('E42').ClearFormats()
$ws.Range('E43').Value = "'  +0.05%  "
$ws.Range('E43').ClearFormats()
$ws.Range('D44').Value = "'101.26"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = "'  +0.53%  "
$ws.Range('E44').ClearFormats()
$ws.Range('D45').Value = "'1.868.14"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = "'  +4.61%  "
$ws.Range('E45').ClearFormats()
$ws.Range('E46').Value = "'  +1.67%  "
$ws.Range('E46').ClearFormats()
$ws.Range('D47').Value = "'59.20"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = "'  +3.35%  "
$ws.Range('E47').ClearFormats()
$ws.Range('D48').Value = "'8.244"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = "'  +3.99%  "
$ws.Range('E48').ClearFormats()
$ws.Range('D49').Value = "'0.4433"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = "'  +2.35%  "
$ws.Range('E49').ClearFormats()
$ws.Range('D50').Value = "'1.006"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = "'  +0.86%  "
$ws.Range('E50').ClearFormats()
$ws.Range('D51').Value = "'0.05250"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = "'  +2.18%  "
$ws.Range('E51').ClearFormats()
